$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.012465150128476
$ws.Range("D2").Value = 1.01951604064521
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.010740916476713
$ws.Range("I2").Value = 1.025470798456751
$ws.Range("J2").Value = 1.01770770969491
$ws.Range("K2").Value = 1.022359057194457
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.01361019273714
$ws.Range("N2").Value = 1.010040949031979
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.013345178671971
$ws.Range("D3").Value = 1.020145884982382
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.012265847099409
$ws.Range("I3").Value = 1.025566134964549
$ws.Range("J3").Value = 1.018221816741267
$ws.Range("K3").Value = 1.022795427714891
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.014937246691773
$ws.Range("N3").Value = 1.010211950226538
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.013914539770927
$ws.Range("D4").Value = 1.020553129343506
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.013252413402638
$ws.Range("I4").Value = 1.025626311157497
$ws.Range("J4").Value = 1.018553805669534
$ws.Range("K4").Value = 1.023076795657036
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.015795299103681
$ws.Range("N4").Value = 1.01032233751154
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.014153880806983
$ws.Range("D5").Value = 1.020724260529732
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.013667131301661
$ws.Range("I5").Value = 1.025651246715274
$ws.Range("J5").Value = 1.018693212677721
$ws.Range("K5").Value = 1.023194844593337
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.016155875392288
$ws.Range("N5").Value = 1.010368681547543
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.014194066139178
$ws.Range("D6").Value = 1.020752989824893
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.01373676238266
$ws.Range("I6").Value = 1.025655412237765
$ws.Range("J6").Value = 1.018716610275939
$ws.Range("K6").Value = 1.023214651551558
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.016216409129126
$ws.Range("N6").Value = 1.010376459231554
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.013917737929911
$ws.Range("D7").Value = 1.020555416299983
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.013257955013004
$ws.Range("I7").Value = 1.02562664577248
$ws.Range("J7").Value = 1.018555669066086
$ws.Range("K7").Value = 1.023078373969009
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.01580011772158
$ws.Range("N7").Value = 1.010322957009592
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.012762575603121
$ws.Range("D8").Value = 1.019728962562579
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.011256311326515
$ws.Range("I8").Value = 1.025503330939751
$ws.Range("J8").Value = 1.017881593815127
$ws.Range("K8").Value = 1.022506735998869
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.014058812006992
$ws.Range("N8").Value = 1.010098793875453
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.010726452973645
$ws.Range("D9").Value = 1.018270330221581
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.007727660642969
$ws.Range("I9").Value = 1.025274461403161
$ws.Range("J9").Value = 1.016688638248549
$ws.Range("K9").Value = 1.021491850308427
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.01098531139976
$ws.Range("N9").Value = 1.009701786810282
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.009368651740673
$ws.Range("D10").Value = 1.017296399074316
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.005373897520623
$ws.Range("I10").Value = 1.025114118738875
$ws.Range("J10").Value = 1.015889874834531
$ws.Range("K10").Value = 1.02081018198698
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.008932610968912
$ws.Range("N10").Value = 1.0094357720568
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.008780616866363
$ws.Range("D11").Value = 1.016874327407501
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.004354304818568
$ws.Range("I11").Value = 1.025042850876049
$ws.Range("J11").Value = 1.015543180544393
$ws.Range("K11").Value = 1.02051381190605
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.008042830255681
$ws.Range("N11").Value = 1.009320266550671
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.008562179664324
$ws.Range("D12").Value = 1.016717498840403
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.003975516668315
$ws.Range("I12").Value = 1.025016102856504
$ws.Range("J12").Value = 1.015414278865622
$ws.Range("K12").Value = 1.020403546359815
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.007712178328774
$ws.Range("N12").Value = 1.009277314688424
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.008609035860887
$ws.Range("D13").Value = 1.016751141458807
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.004056771101
$ws.Range("I13").Value = 1.025021852881515
$ws.Range("J13").Value = 1.015441934325729
$ws.Range("K13").Value = 1.020427206857729
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.007783111061695
$ws.Range("N13").Value = 1.009286530179553
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.008762561083269
$ws.Range("D14").Value = 1.016861364973373
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.004322995455115
$ws.Range("I14").Value = 1.025040645504566
$ws.Range("J14").Value = 1.015532528021331
$ws.Range("K14").Value = 1.020504701002291
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.008015501509361
$ws.Range("N14").Value = 1.009316717114977
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.008857151065145
$ws.Range("D15").Value = 1.016929270400981
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.004487016129489
$ws.Range("I15").Value = 1.025052187699379
$ws.Range("J15").Value = 1.0155883293546
$ws.Range("K15").Value = 1.02055242379186
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.008158665195096
$ws.Range("N15").Value = 1.009335309925249
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.009407675601918
$ws.Range("D16").Value = 1.017324403201073
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.005441555616082
$ws.Range("I16").Value = 1.025118809839256
$ws.Range("J16").Value = 1.01591286641251
$ws.Range("K16").Value = 1.020829825765938
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.008991642358061
$ws.Range("N16").Value = 1.009443431048831
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.009752978765827
$ws.Range("D17").Value = 1.017572165532062
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.006040202117356
$ws.Range("I17").Value = 1.025160108118636
$ws.Range("J17").Value = 1.016116219123302
$ws.Range("K17").Value = 1.021003510754397
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.00951388938703
$ws.Range("N17").Value = 1.009511167076114
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.00995437881588
$ws.Range("D18").Value = 1.017716646998083
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.006389344539216
$ws.Range("I18").Value = 1.025184019308314
$ws.Range("J18").Value = 1.016234751778809
$ws.Range("K18").Value = 1.021104702161031
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.009818416146253
$ws.Range("N18").Value = 1.009550645540729
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.010023049430896
$ws.Range("D19").Value = 1.01776590563764
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.006508386835015
$ws.Range("I19").Value = 1.025192142308913
$ws.Range("J19").Value = 1.016275154863213
$ws.Range("K19").Value = 1.021139186113461
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.009922236612051
$ws.Range("N19").Value = 1.009564101450109
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.009715931981942
$ws.Range("D20").Value = 1.01754558650184
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.005975976983648
$ws.Range("I20").Value = 1.025155695551425
$ws.Range("J20").Value = 1.01609440952989
$ws.Range("K20").Value = 1.020984888006523
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.009457866683907
$ws.Range("N20").Value = 1.009503902828596
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.008717352143424
$ws.Range("D21").Value = 1.016828908335452
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.004244600891337
$ws.Range("I21").Value = 1.025035119161451
$ws.Range("J21").Value = 1.015505853859475
$ws.Range("K21").Value = 1.020481885888867
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.007947072436433
$ws.Range("N21").Value = 1.009307829139176
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.00808941926877
$ws.Range("D22").Value = 1.016378001735326
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.003155630909402
$ws.Range("I22").Value = 1.024957711485288
$ws.Range("J22").Value = 1.015135088159735
$ws.Range("K22").Value = 1.020164584288001
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.006996318615494
$ws.Range("N22").Value = 1.009184272264412
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.008422306578579
$ws.Range("D23").Value = 1.016617064311926
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.003732953106246
$ws.Range("I23").Value = 1.024998898010285
$ws.Range("J23").Value = 1.01533170611373
$ws.Range("K23").Value = 1.020332890762525
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.007500414278783
$ws.Range("N23").Value = 1.009249798410093
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.009732671857927
$ws.Range("D24").Value = 1.017557596526902
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.006004997672832
$ws.Range("I24").Value = 1.025157689948729
$ws.Range("J24").Value = 1.016104264592001
$ws.Range("K24").Value = 1.020993303183371
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.009483181212172
$ws.Range("N24").Value = 1.009507185324536
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.011252907743381
$ws.Range("D25").Value = 1.01864769112668
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.008640112080344
$ws.Range("I25").Value = 1.025334999255263
$ws.Range("J25").Value = 1.016997655778581
$ws.Range("K25").Value = 1.021755119241601
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.011780517073646
$ws.Range("N25").Value = 1.009804659664095
